# Apply the edits described by the diff:
# - Set B3=16, B4=10, B5=10 on the "sell" sheet
# - Select/activate cell B8 on the "sell" sheet
# - Make the "sell" sheet the active tab (instead of "回收")

$wb = $excel.ActiveWorkbook

$sell = $wb.Worksheets.Item("sell")

# Fill in the new quantity values in column B
$sell.Range("B3").Value = 16
$sell.Range("B4").Value = 10
$sell.Range("B5").Value = 10

# Activate the "sell" sheet and select B8 on it, which becomes the
# active tab / active cell for the workbook
$sell.Activate()
$sell.Range("B8").Select()
